# Swap the "category-code" (column F) and "category-name" (column G)
# values for every row on the active sheet, including the header row.
#
# The source data had these two columns reversed (category-code before
# category-name); this corrects the column order by exchanging the
# values that live in F and G, row by row, leaving every other column
# (A-E) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row   # -4162 = xlUp
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $fValue = $fCell.Value()
    $gValue = $gCell.Value()

    $fCell.Value = $gValue
    $gCell.Value = $fValue
}
